$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(42606.52685185185,   16, 75, 23, 25, 75, 1882, 3360, 427, 73, 23, 2, 6, "Bag"),
    @(42606.542592592596,  14, 75, 24, 25, 75, 1568, 3367, 427, 73, 24, 2, 6, "Bag"),
    @(42606.551504629628,  14, 75, 24, 25, 75, 1695, 3367, 427, 73, 24, 2, 6, "Bag"),
    @(42606.554189814815,  14, 75, 24, 25, 75, 1592, 3367, 427, 73, 24, 2, 6, "Bag"),
    @(42606.555902777778,  14, 75, 24, 25, 75, 1629, 3369, 427, 73, 24, 2, 6, "Bag"),
    @(42606.55945601852,   14, 75, 24, 25, 75, 1632, 3367, 427, 73, 24, 2, 6, "Bag"),
    @(42606.566932870373,  14, 75, 24, 25, 75, 1522, 3367, 427, 73, 24, 2, 6, "Bag"),
    @(42606.571203703701,  14, 75, 24, 25, 75, 1560, 3367, 427, 73, 24, 2, 6, "Bag")
)

$row = 3
foreach ($rec in $data) {
    $col = 1
    foreach ($val in $rec) {
        $ws.Cells.Item($row, $col).Value = $val
        $col++
    }
    $row++
}
